# Update the handoff/handback timestamps in the zh-cn and de-de report sheets
# to reflect a newly generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 08:52:42"
$wsZhCn.Range("H2").Value = "2016-03-24 08:53:07"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 08:52:46"
$wsDeDe.Range("H2").Value = "2016-03-24 08:53:14"
